$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4: enter the Session 3 (D&C) mark for the student
$ws.Range("D4").Value = 10

# D5: add feedback comment for the Session 3 (D&C) column (merged D5:D12)
$ws.Range("D5").Value = "Quadratic times are a little bit weird but you did a good work"

# Update the active selection to match the edited feedback cell C5 (merged C5:C12)
$ws.Range("C5:C12").Select()
